$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 502
